$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are kept as text, preserving exact formatting
# (leading/trailing zeros, thousands separators as dots, etc.)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "63.916.65"
$ws.Range("E2").Value2 = "  +0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.137.75"
$ws.Range("E3").Value2 = "  +0.64%  "

$ws.Range("E4").Value2 = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "589.71"
$ws.Range("E5").Value2 = "  +0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "145.32"
$ws.Range("E6").Value2 = "  -0.59%  "

$ws.Range("E7").Value2 = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.130.93"
$ws.Range("E8").Value2 = "  +0.65%  "

$ws.Range("E9").Value2 = "  -0.25%  "

$ws.Range("E10").Value2 = "  -1.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "5.90"
$ws.Range("E11").Value2 = "  +2.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.458"
$ws.Range("E12").Value2 = "  -1.65%  "

$ws.Range("E13").Value2 = "  -2.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "37.29"
$ws.Range("E14").Value2 = "  +0.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.658.00"
$ws.Range("E15").Value2 = "  +0.66%  "

$ws.Range("E16").Value2 = "  -1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "7.33"
$ws.Range("E17").Value2 = "  +2.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "63.763.71"
$ws.Range("E18").Value2 = "  +0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "3.130.33"
$ws.Range("E19").Value2 = "  +0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "468.23"
$ws.Range("E20").Value2 = "  +0.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "14.35"

$ws.Range("E22").Value2 = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.54"
$ws.Range("E23").Value2 = "  +0.08%  "

$ws.Range("B24").Value2 = "Litecoin"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "81.62"
$ws.Range("E24").Value2 = "  -0.50%  "

$ws.Range("B25").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "12.97"
$ws.Range("E25").Value2 = "  -1.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.32"
$ws.Range("E26").Value2 = "  +7.10%  "

$ws.Range("E27").Value2 = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "9.84"
$ws.Range("E28").Value2 = "  +9.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.46"
$ws.Range("E29").Value2 = "  +8.67%  "

$ws.Range("B30").Value2 = "ImmutableX"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "2.24"
$ws.Range("E30").Value2 = "  +0.60%  "

$ws.Range("B31").Value2 = "PancakeSwap"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "2.71"
$ws.Range("E31").Value2 = "  +0.35%  "

$ws.Range("E32").Value2 = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "27.72"
$ws.Range("E33").Value2 = "  +2.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.110"
$ws.Range("E34").Value2 = "  +0.64%  "

$ws.Range("E35").Value2 = "  -3.57%  "

$ws.Range("E36").Value2 = "  +1.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "6.15"
$ws.Range("E37").Value2 = "  +1.18%  "

$ws.Range("E38").Value2 = "  -2.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "3.21"
$ws.Range("E39").Value2 = "  -5.80%  "

$ws.Range("B40").Value2 = "Cosmos"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "9.36"
$ws.Range("E40").Value2 = "  +7.72%  "

$ws.Range("B41").Value2 = "OKB"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "51.21"
$ws.Range("E41").Value2 = "  +0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "454.48"
$ws.Range("E42").Value2 = "  +1.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.292"
$ws.Range("E43").Value2 = "  +5.60%  "

$ws.Range("E44").Value2 = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "2.917.80"
$ws.Range("E45").Value2 = "  +1.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "40.00"
$ws.Range("E46").Value2 = "  +11.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.107"
$ws.Range("E47").Value2 = "  -2.80%  "

$ws.Range("E48").Value2 = "  +6.78%  "

$ws.Range("E50").Value2 = "  +2.81%  "

$ws.Range("E51").Value2 = "  -0.63%  "

